# Update cryptocurrency price (D) and 1h-volume-change (E) columns
# to the latest scraped figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.191.01'
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").Value = '2.056.91'
$ws.Range("E3").Value = '  -0.73%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'" + '248.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.00%  '
$ws.Range("E6").Value = '  -1.71%  '
$ws.Range("D7").Value = "'" + '58.09'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -6.51%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -2.22%  '
$ws.Range("D10").Value = "'" + '0.0783'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.76%  '
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("D12").Value = "'" + '15.98'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.82%  '
$ws.Range("D13").Value = '2.355.22'
$ws.Range("E13").Value = '  -0.83%  '
$ws.Range("D14").Value = "'" + '0.840'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.78%  '
$ws.Range("D15").Value = "'" + '5.71'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.89%  '
$ws.Range("D16").Value = '2.056.92'
$ws.Range("E16").Value = '  -0.76%  '
$ws.Range("D17").Value = "'" + '18.15'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +17.01%  '
$ws.Range("D18").Value = '37.213.07'
$ws.Range("E18").Value = '  +0.34%  '
$ws.Range("D19").Value = "'" + '74.88'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").Value = '0.0₃0900'
$ws.Range("E20").Value = '  -3.06%  '
$ws.Range("D21").Value = "'" + '5.38'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.24%  '
$ws.Range("D22").Value = "'" + '237.78'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.23%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").Value = "'" + '2.47'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.21%  '
$ws.Range("E25").Value = '  -8.58%  '
$ws.Range("D26").Value = "'" + '169.65'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.14%  '
$ws.Range("D27").Value = "'" + '9.41'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.59%  '
$ws.Range("D28").Value = "'" + '20.11'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.53%  '
$ws.Range("E29").Value = '  -1.73%  '
$ws.Range("D30").Value = "'" + '1.13'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.65%  '
$ws.Range("E31").Value = '  -1.22%  '
$ws.Range("D32").Value = "'" + '0.0619'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.15%  '
$ws.Range("D33").Value = "'" + '4.51'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.95%  '
$ws.Range("D34").Value = "'" + '0.0907'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.88%  '
$ws.Range("E35").Value = '  +0.14%  '
$ws.Range("D36").Value = "'" + '2.27'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.89%  '
$ws.Range("D37").Value = "'" + '1.76'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.40%  '
$ws.Range("E38").Value = '  -2.68%  '
$ws.Range("D39").Value = "'" + '3.18'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +12.98%  '
$ws.Range("E40").Value = '  -10.02%  '
$ws.Range("E41").Value = '  +15.68%  '
$ws.Range("D42").Value = "'" + '0.0223'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.51%  '
$ws.Range("D43").Value = "'" + '1.15'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.88%  '
$ws.Range("D44").Value = "'" + '17.19'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.34%  '
$ws.Range("D45").Value = "'" + '96.28'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.11%  '
$ws.Range("E46").Value = '  -1.84%  '
$ws.Range("D47").Value = "'" + '2.92'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.45%  '
$ws.Range("D48").Value = '1.279.21'
$ws.Range("E48").Value = '  -2.28%  '
$ws.Range("D49").Value = "'" + '6.84'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.63%  '
$ws.Range("D50").Value = '2.240.46'
$ws.Range("E50").Value = '  -0.23%  '
$ws.Range("D51").Value = "'" + '43.93'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.54%  '
